$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (from 2023-09-09 / 45178 to 2023-09-10 / 45179) for every data row (2-125).
$oldSerial = 45178
$newSerial = 45179

for ($row = 2; $row -le 125; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq $oldSerial) {
        $cell.Value2 = $newSerial
    }
}
